# Insert two new rows at 1080-1081, shifting all existing rows (old 1080..1179)
# down to (new 1082..1181). This also grows the sheet dimension from
# A1:R1179 to A1:R1181 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1080:A1081").EntireRow.Insert()

# --- New row 1080: Tomate, Larga vida, Primera, fecha 45166 ---
$ws.Range("A1080").Value2 = 5
$ws.Range("B1080").Value2 = "Macroferia Regional de Talca"
$ws.Range("C1080").Value2 = "Maule"
$ws.Range("D1080").Value2 = 45166
$ws.Range("D1080").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1080").Value2 = 7
$ws.Range("F1080").Value2 = 100112020
$ws.Range("G1080").Value2 = "Tomate"
$ws.Range("H1080").Value2 = "Larga vida"
$ws.Range("I1080").Value2 = "Primera"
$ws.Range("J1080").Value2 = 2000
$ws.Range("K1080").Value2 = 13000
$ws.Range("L1080").Value2 = 13000
$ws.Range("M1080").Value2 = 13000
$ws.Range("N1080").Value2 = "$/bandeja 18 kilos"
$ws.Range("O1080").Value2 = "Región de Arica y Parinacota"
$ws.Range("P1080").Value2 = 722
$ws.Range("Q1080").Value2 = 18
$ws.Range("R1080").Value2 = "Hortaliza"

# --- New row 1081: Tomate, Larga vida, Segunda, fecha 45166 ---
$ws.Range("A1081").Value2 = 5
$ws.Range("B1081").Value2 = "Macroferia Regional de Talca"
$ws.Range("C1081").Value2 = "Maule"
$ws.Range("D1081").Value2 = 45166
$ws.Range("D1081").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1081").Value2 = 7
$ws.Range("F1081").Value2 = 100112020
$ws.Range("G1081").Value2 = "Tomate"
$ws.Range("H1081").Value2 = "Larga vida"
$ws.Range("I1081").Value2 = "Segunda"
$ws.Range("J1081").Value2 = 800
$ws.Range("K1081").Value2 = 10000
$ws.Range("L1081").Value2 = 10000
$ws.Range("M1081").Value2 = 10000
$ws.Range("N1081").Value2 = "$/bandeja 18 kilos"
$ws.Range("O1081").Value2 = "Región de Arica y Parinacota"
$ws.Range("P1081").Value2 = 556
$ws.Range("Q1081").Value2 = 18
$ws.Range("R1081").Value2 = "Hortaliza"
